$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2432.697
$ws.Range("I112").Value = 995
$ws.Range("J112").Value = 2525.4517
$ws.Range("K112").Value = 2985
$ws.Range("L112").Value = 7576.355100000001
$ws.Range("M112").Value = -1877
$ws.Range("N112").Value = -9792.355100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 702.9091
$ws.Range("I129").Value = 273
$ws.Range("J129").Value = 948.5714
$ws.Range("K129").Value = 819
$ws.Range("L129").Value = 2845.7142
$ws.Range("M129").Value = 4181
$ws.Range("N129").Value = -12845.7142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7999
$ws.Range("I141").Value = 395
$ws.Range("K141").Value = 1185
$ws.Range("M141").Value = 3995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1227186.2
$ws.Range("I2").Value = 1519.6666
$ws.Range("J2").Value = 3269964
$ws.Range("K2").Value = 1519.6666
$ws.Range("L2").Value = 3269964
$ws.Range("M2").Value = -1406.6666
$ws.Range("N2").Value = -3270190

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10600.143
$ws.Range("I28").Value = 880
$ws.Range("J28").Value = 34900.5
$ws.Range("K28").Value = 880
$ws.Range("L28").Value = 34900.5
$ws.Range("M28").Value = -688
$ws.Range("N28").Value = -35284.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9056.799999999999
$ws.Range("I32").Value = 7500.9316
$ws.Range("J32").Value = 20466.5
$ws.Range("K32").Value = 7500.9316
$ws.Range("L32").Value = 20466.5
$ws.Range("M32").Value = -7213.9316
$ws.Range("N32").Value = -21040.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 10600.143
$ws.Range("I99").Value = 880
$ws.Range("J99").Value = 34900.5
$ws.Range("K99").Value = 880
$ws.Range("L99").Value = 34900.5
$ws.Range("M99").Value = 2115
$ws.Range("N99").Value = -40890.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1227186.2
$ws.Range("I116").Value = 1519.6666
$ws.Range("J116").Value = 3269964
$ws.Range("K116").Value = 1519.6666
$ws.Range("L116").Value = 3269964
$ws.Range("M116").Value = 774.3334
$ws.Range("N116").Value = -3274552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 26100
$ws.Range("J117").Value = 26100
$ws.Range("L117").Value = 26100
$ws.Range("N117").Value = -35278

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H129").Value = 37999.6
$ws.Range("J129").Value = 37999.6
$ws.Range("L129").Value = 37999.6
$ws.Range("N129").Value = -47999.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1227186.2
$ws.Range("I3").Value = 1519.6666
$ws.Range("J3").Value = 3269964
$ws.Range("K3").Value = 1519.6666
$ws.Range("L3").Value = 3269964
$ws.Range("M3").Value = -1405.6666
$ws.Range("N3").Value = -3270192

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 194
$ws.Range("I80").Value = 126.333336
$ws.Range("J80").Value = 244.75
$ws.Range("K80").Value = 126.333336
$ws.Range("L80").Value = 244.75
$ws.Range("M80").Value = 871.666664
$ws.Range("N80").Value = -2240.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 194
$ws.Range("I83").Value = 126.333336
$ws.Range("J83").Value = 244.75
$ws.Range("K83").Value = 631.66668
$ws.Range("L83").Value = 1223.75
$ws.Range("M83").Value = 4360.33332
$ws.Range("N83").Value = -11207.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1283146.1
$ws.Range("I107").Value = 1411062.8
$ws.Range("J107").Value = 3980
$ws.Range("K107").Value = 1411062.8
$ws.Range("L107").Value = 3980
$ws.Range("M107").Value = -1409142.8
$ws.Range("N107").Value = -7820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 53437.25
$ws.Range("J129").Value = 53437.25
$ws.Range("L129").Value = 53437.25
$ws.Range("N129").Value = -63437.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 46090.668
$ws.Range("J139").Value = 46090.668
$ws.Range("L139").Value = 46090.668
$ws.Range("N139").Value = -56370.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10021050
$ws.Range("I6").Value = 16701367
$ws.Range("J6").Value = 575
$ws.Range("K6").Value = 16701367
$ws.Range("L6").Value = 575
$ws.Range("M6").Value = -16701254
$ws.Range("N6").Value = -801

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3534.4255
$ws.Range("I58").Value = 1275.3572
$ws.Range("J58").Value = 4492.8184
$ws.Range("K58").Value = 1275.3572
$ws.Range("L58").Value = 4492.8184
$ws.Range("M58").Value = -1072.3572
$ws.Range("N58").Value = -4898.8184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3534.4255
$ws.Range("I136").Value = 1275.3572
$ws.Range("J136").Value = 4492.8184
$ws.Range("K136").Value = 3826.0716
$ws.Range("L136").Value = 13478.4552
$ws.Range("M136").Value = -1276.0716
$ws.Range("N136").Value = -18578.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2247.7646
$ws.Range("I109").Value = 1243.4
$ws.Range("J109").Value = 2666.25
$ws.Range("K109").Value = 3730.2
$ws.Range("L109").Value = 7998.75
$ws.Range("M109").Value = -2690.2
$ws.Range("N109").Value = -10078.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 3009
$ws.Range("I110").Value = 2013.5
$ws.Range("K110").Value = 6040.5
$ws.Range("M110").Value = -1950.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 41668864
$ws.Range("I112").Value = 2090
$ws.Range("J112").Value = 55557790
$ws.Range("K112").Value = 6270
$ws.Range("L112").Value = 166673370
$ws.Range("M112").Value = -5162
$ws.Range("N112").Value = -166675586

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -23744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 7655
$ws.Range("J117").Value = 7655
$ws.Range("L117").Value = 7655
$ws.Range("N117").Value = -14539

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 20000
$ws.Range("J131").Value = 20000
$ws.Range("L131").Value = 20000
$ws.Range("N131").Value = -30080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 32533.334
$ws.Range("J118").Value = 32533.334
$ws.Range("L118").Value = 32533.334
$ws.Range("N118").Value = -35847.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 54003.332
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 54003.332
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 54003.332
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -54339.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 14998.667
$ws.Range("J39").Value = 14998.667
$ws.Range("L39").Value = 14998.667
$ws.Range("N39").Value = -15824.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 727.36365
$ws.Range("I113").Value = 557.4286
$ws.Range("J113").Value = 1024.75
$ws.Range("K113").Value = 1672.2858
$ws.Range("L113").Value = 3074.25
$ws.Range("M113").Value = 497.7142000000001
$ws.Range("N113").Value = -7414.25
